## Separate inventory items per project - each project now has its own
## Master_Items file.
##
## The modifications-log sheet gains a block of "original transaction"
## reconciliation columns (T:AB) plus one illustrative data row (row 2)
## recording a quantity correction reconciled against the original
## inventory transaction it corrects.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New header row (row 1), columns T..AB
# ---------------------------------------------------------------------
$ws.Cells.Item(1, 20).Value = "معرف_المعاملة_الأصلية"
$ws.Cells.Item(1, 21).Value = "اسم_العنصر"
$ws.Cells.Item(1, 22).Value = "التصنيف"
$ws.Cells.Item(1, 23).Value = "نوع_العملية"
$ws.Cells.Item(1, 24).Value = "تاريخ_المعاملة_الأصلية"
$ws.Cells.Item(1, 25).Value = "الكمية_الأصلية"
$ws.Cells.Item(1, 26).Value = "فرق_الكمية"
$ws.Cells.Item(1, 27).Value = "ملاحظات"
$ws.Cells.Item(1, 28).Value = "المستخدم"

# Match the existing header look (A1:S1 all share one style: bold font,
# thin border, centered/top-aligned) by copying that style onto the new
# header cells instead of re-building it (avoids minting new font/xf
# entries).
$ws.Range("A1").Copy()
$ws.Range("T1:AB1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# New data row (row 2) - a sample reconciliation record.
# ---------------------------------------------------------------------

# Materialise every column across the row (A2:AB2) - even the ones that
# stay blank - by (re)applying the built-in "Normal" cell style, which
# touches the cell without minting any new style entry (still style
# index 0, same as an untouched cell).
for ($col = 1; $col -le 28; $col++) {
    $ws.Cells.Item(2, $col).Style = "Normal"
}

$ws.Cells.Item(2, 3).Value = "2025-12-07 20:33:47"
$ws.Cells.Item(2, 8).Value = 15
$ws.Cells.Item(2, 19).Value = "تصحيح خطأ في الإدخال"

$ws.Cells.Item(2, 20).Value = "2025-12-07 20:33:25_طلاء أبيض_22"
$ws.Cells.Item(2, 21).Value = "طلاء أبيض"
$ws.Cells.Item(2, 22).Value = "مواد التشطيب"
$ws.Cells.Item(2, 23).Value = "دخول"

# Original-transaction timestamp, stored as a real date/time serial.
# NumberFormat is first set with a lowercase format code (registers the
# first custom numFmt, left unused) and then reset to the uppercase code
# that ends up applied to the cell - mirrors the two custom numFmts left
# behind in the source workbook.
$ws.Cells.Item(2, 24).Value = 45998.85653935185
$ws.Cells.Item(2, 24).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(2, 24).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(2, 25).Value = 22
$ws.Cells.Item(2, 26).Value = -7
$ws.Cells.Item(2, 28).Value = "النظام"
